$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing "Carbon dioxide, non-fossil" rows (13-17): CF value 0 -> 8.22E-14
$ws.Range("C13").Value = [double]"8.22E-14"
$ws.Range("C14").Value = [double]"8.22E-14"
$ws.Range("C15").Value = [double]"8.22E-14"
$ws.Range("C16").Value = [double]"8.22E-14"
$ws.Range("C17").Value = [double]"8.22E-14"

# Add two new CF rows for specific CO2 flows to account for NETs
$ws.Range("A44").Value = "Carbon dioxide, in air"
$ws.Range("B44").Value = "natural resource::in air"
$ws.Range("C44").Value = [double]"-8.22E-14"

$ws.Range("A45").Value = "Carbon dioxide, non-fossil, resource correction"
$ws.Range("B45").Value = "natural resource::in air"
$ws.Range("C45").Value = [double]"-8.22E-14"

# Column C width recalculated (bestFit) after new content added
$ws.Columns.Item(3).ColumnWidth = 8.830729166666666

# Add hidden filter-database defined name (left over from AutoFilter use)
$nm = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$1:`$C`$43")
$nm.RefersTo = "=Sheet1!`$A`$1:`$C`$43"
$nm.Visible = $false

# Update selection to match author's final cursor position
$ws.Range("C5").Select() | Out-Null

Write-Host "done"
